$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data row that was incorrectly left down at row 9 (with a big empty gap
# under the header) needs to sit directly under the header row instead.
# Move it (row 9 -> row 2); rows 10/11 stay exactly where they are.
$ws.Range("A9:K9").Copy($ws.Range("A2"))
$ws.Range("A9:K9").ClearContents()

# Now that real data lives right under the header, give the columns that
# hold it a sensible best-fit width instead of the default.
$ws.Columns("B:G").AutoFit()

$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 33.333333333333336
$ws.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws.Columns.Item(7).ColumnWidth = 9.666666666666666
